$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the very top; everything shifts down by one.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Article"
$ws.Range("B1").Value = "Annotation"

# Re-label the article cells in place (drop the old "Article N - " prefixes),
# touching them in this order: Article 5, then Article 9.3, then Article 8.4.
$ws.Range("A5").Value = "5. Udemy’s Rights to Content You Post"
$ws.Range("A9").Value = "9.3 Limitation of liability"
$ws.Range("A7").Value = "8.4 Payments and billing"

# Remove the old ellipsis placeholder rows (now at 8, 6, 4, 3, 2) bottom-up
# so row numbers of the rows not-yet-deleted stay stable.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

# Insert a second annotation row right under "8.4 Payments and billing".
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "8.4 Payments and billing"
$ws.Range("B4").Value = "Subscription plan can only be changed by user"

# Column widths for the new two-column layout (47 and 77.85546875 "characters").
$ws.Columns.Item(1).ColumnWidth = 46.16666666666667
$ws.Columns.Item(2).ColumnWidth = 77

$ws.Range("B3").Select()
